$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared-string label for the macrophage cluster (Resolving-Mac -> Inflammatory-Mac).
# Rows 4/8/12 move from "MuSCs" to the macrophage cluster, rows 5/9/13 move the other way,
# matching the corrected target-cluster assignment that came with the new TPM values below.

# New TPM-derived values for the Ntf3-Ngfr sheet (columns D, G:T), rows 2-13
$rowData = @(
    @{ Row = 2; D = "ECs"; Vals = @(1.679513333333333, 5.038539999999999, 0.376631045782902, 0.3928053077137587, 2, 0.6666666666666666, 0.08108666666666665, 0.24326, 0.02056149724823249, 0.02939174694341418, 0.1361861378222222, 1.2256752404, 0.007744098211464063, 0.01154523420235273) },
    @{ Row = 3; D = "FAPs"; Vals = @(1.679513333333333, 5.038539999999999, 0.376631045782902, 0.3928053077137587, 2, 0.6666666666666666, 0.090517, 0.271551, 0.02295278771378271, 0.03280999043916413, 0.1520245083933333, 1.36822057554, 0.008644732440274925, 0.01288793839054135) },
    @{ Row = 4; D = "Inflammatory-Mac"; Vals = @(1.679513333333333, 5.038539999999999, 0.376631045782902, 0.3928053077137587, 2, 0.6666666666666666, 0.2176356666666667, 0.652907, 0.05518681856388939, 0.07888710565478799, 0.3655220039755555, 3.28969803578, 0.02078506918914893, 0.03098727381137679) },
    @{ Row = 5; D = "MuSCs"; Vals = @(1.679513333333333, 5.038539999999999, 0.376631045782902, 0.3928053077137587, 2, 1, 3.5543775, 7.108755, 0.9012988964740954, 0.8589111569626339, 5.96962440295, 35.8177464177, 0.339457145942014, 0.3373848613094879) },
    @{ Row = 6; D = "ECs"; Vals = @(2.228940333333334, 6.686821, 0.4998401096732527, 0.5213055330575571, 2, 0.6666666666666666, 0.08108666666666665, 0.24326, 0.02056149724823249, 0.02939174694341418, 0.1807373418288889, 1.62663607646, 0.01027746103960281, 0.01532208030782935) },
    @{ Row = 7; D = "FAPs"; Vals = @(2.228940333333334, 6.686821, 0.4998401096732527, 0.5213055330575571, 2, 0.6666666666666666, 0.090517, 0.271551, 0.02295278771378271, 0.03280999043916413, 0.2017569921523334, 1.815812929371, 0.01147272392816403, 0.01710402955550181) },
    @{ Row = 8; D = "Inflammatory-Mac"; Vals = @(2.228940333333334, 6.686821, 0.4998401096732527, 0.5213055330575571, 2, 0.6666666666666666, 0.2176356666666667, 0.652907, 0.05518681856388939, 0.07888710565478799, 0.4850969154052223, 4.365872238647, 0.02758458544349237, 0.04112428466473708) },
    @{ Row = 9; D = "MuSCs"; Vals = @(2.228940333333334, 6.686821, 0.4998401096732527, 0.5213055330575571, 2, 1, 3.5543775, 7.108755, 0.9012988964740954, 0.8589111569626339, 7.922495369642501, 47.53497221785501, 0.4505053392619935, 0.4477551385294889) },
    @{ Row = 10; D = "ECs"; Vals = @(0.550853, 1.101706, 0.1235288445438454, 0.08588915922868416, 2, 0.6666666666666666, 0.08108666666666665, 0.24326, 0.02056149724823249, 0.02939174694341418, 0.04466683359333333, 0.26800100156, 0.002539937997165616, 0.002524432433232091) },
    @{ Row = 11; D = "FAPs"; Vals = @(0.550853, 1.101706, 0.1235288445438454, 0.08588915922868416, 2, 0.6666666666666666, 0.090517, 0.271551, 0.02295278771378271, 0.03280999043916413, 0.049861561001, 0.299169366006, 0.002835331345343748, 0.002818022493120973) },
    @{ Row = 12; D = "Inflammatory-Mac"; Vals = @(0.550853, 1.101706, 0.1235288445438454, 0.08588915922868416, 2, 0.6666666666666666, 0.2176356666666667, 0.652907, 0.05518681856388939, 0.07888710565478799, 0.1198852598903333, 0.7193115593420001, 0.006817163931248094, 0.006775547178674116) },
    @{ Row = 13; D = "MuSCs"; Vals = @(0.550853, 1.101706, 0.1235288445438454, 0.08588915922868416, 2, 1, 3.5543775, 7.108755, 0.9012988964740954, 0.8589111569626339, 1.9579395090075, 7.831758036030001, 0.1113364112700879, 0.07377115712365699) }
)

foreach ($entry in $rowData) {
    $r = $entry.Row
    $ws.Cells.Item($r, 4).Value = $entry.D
    for ($i = 0; $i -lt $entry.Vals.Length; $i++) {
        $ws.Cells.Item($r, 7 + $i).Value = $entry.Vals[$i]
    }
}

